$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert M2:M4 from text to numeric (large id -> double), per diff ---
$ws.Range("M2").Value = 54926309110740000
$ws.Range("M3").Value = 55224007200130201
$ws.Range("M4").Value = 55224007200130301

# --- Add rows 5-10 to the sheet ---

# Columns A (PEDIDO) and M (INSTALACION) hold long, all-digit values that
# must stay text (leading zeros, > 15 significant digits). Mark them as
# Text before assigning so Excel does not coerce them into numbers, then
# drop back to the Normal cell style (the text flag on the cell persists).
$ws.Range("A5:A10").NumberFormat = "@"
$ws.Range("M5:M10").NumberFormat = "@"

# F (FECHA_INICIO_ANS) and S (FECHA_LIMITE_ANS) use the same date/time format as the existing rows
$ws.Range("F5:F10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S5:S10").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 5
$ws.Range("A5").Value = "23087278"
$ws.Range("B5").Value = "ENERES"
$ws.Range("C5").Value = "NUEVO"
$ws.Range("D5").Value = "ENERES"
$ws.Range("E5").Value = "18/04/2024 09:10"
$ws.Range("F5").Value = 45939.74375
$ws.Range("G5").Value = 71603529
$ws.Range("H5").Value = "GERMAN DE JESUS MARIN HENAO"
$ws.Range("I5").Value = 4657981
$ws.Range("J5").Value = 3207426955
$ws.Range("K5").Value = "CR 24 CL 57 B -13 (INTERIOR 302 )"
$ws.Range("L5").Value = "MEDELLÍN"
$ws.Range("M5").Value = "055224007200130302"
$ws.Range("N5").Value = "Medellín"
$ws.Range("O5").Value = "ALEGA"
$ws.Range("P5").Value = "METROSUR Vin. Leg. Ref. Concentrada"
$ws.Range("Q5").Value = "Urbano"
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 45951.74375
$ws.Range("T5").Value = "15 días 17:51"
$ws.Range("U5").Value = "VENCIDO"
$ws.Range("V5").Value = "VENCIDO"
$ws.Range("W5").Value = "CERRADO"
$ws.Range("X5").Value = "Ejecutado en Campo"

# Row 6
$ws.Range("A6").Value = "23153422"
$ws.Range("B6").Value = "ENERES"
$ws.Range("C6").Value = "NUEVO"
$ws.Range("D6").Value = "ENEDOM"
$ws.Range("E6").Value = "10/07/2024 13:05"
$ws.Range("F6").Value = 45944.56041666667
$ws.Range("G6").Value = 1035861667
$ws.Range("H6").Value = "LIZANA PATRICIA BEDOYA MEJIA"
$ws.Range("I6").Value = "SIN DATOS"
$ws.Range("J6").Value = 3246410113
$ws.Range("K6").Value = "CR 1 ESTE CL 47 C -18"
$ws.Range("L6").Value = "MEDELLÍN"
$ws.Range("M6").Value = "704021007300180000"
$ws.Range("N6").Value = "Medellín"
$ws.Range("O6").Value = "ARTER"
$ws.Range("P6").Value = "Habilitación Viviendas Metrosur"
$ws.Range("Q6").Value = "Urbano"
$ws.Range("R6").Value = 5
$ws.Range("S6").Value = 45951.56041666667
$ws.Range("T6").Value = "13 días 13:27"
$ws.Range("U6").Value = "VENCIDO"
$ws.Range("V6").Value = "VENCIDO"
$ws.Range("W6").Value = "CERRADO"
$ws.Range("X6").Value = "Ejecutado en Campo"

# Row 7
$ws.Range("A7").Value = "23224950"
$ws.Range("B7").Value = "ENECNX"
$ws.Range("C7").Value = "NUEVO"
$ws.Range("D7").Value = "ENECNX"
$ws.Range("E7").Value = "09/10/2025 10:40"
$ws.Range("F7").Value = 45939.44444444445
$ws.Range("G7").Value = 901761829
$ws.Range("H7").Value = "TAHOE LAQUE SAS TAHOE LAQUE SAS"
$ws.Range("I7").Value = "SIN DATOS"
$ws.Range("J7").Value = 3103589945
$ws.Range("K7").Value = "RURAL_114003250000000000_VEREDA_EL JARDIN"
$ws.Range("L7").Value = "MEDELLÍN"
$ws.Range("M7").Value = "114003250000000000"
$ws.Range("N7").Value = "Medellín"
$ws.Range("O7").Value = "ACREV"
$ws.Range("P7").Value = "Revisor Puntos de Conexión Metrosur"
$ws.Range("Q7").Value = "Rural"
$ws.Range("R7").Value = 4
$ws.Range("S7").Value = 45946.44444444445
$ws.Range("T7").Value = "15 días 10:40"
$ws.Range("U7").Value = "VENCIDO"
$ws.Range("V7").Value = "VENCIDO"
$ws.Range("W7").Value = "CERRADO"
$ws.Range("X7").Value = "Ejecutado en Campo"

# Row 8
$ws.Range("A8").Value = "23252866"
$ws.Range("B8").Value = "ENERES"
$ws.Range("C8").Value = "NUEVO"
$ws.Range("D8").Value = "ENEDOM"
$ws.Range("E8").Value = "08/10/2025 15:54"
$ws.Range("F8").Value = 45938.6625
$ws.Range("G8").Value = 98584619
$ws.Range("H8").Value = "CALEB OBED RAMIREZ MUÑOZ"
$ws.Range("I8").Value = "SIN DATOS"
$ws.Range("J8").Value = 3002300945
$ws.Range("K8").Value = "CR 34 E CL 31 -190 (INTERIOR 119 )"
$ws.Range("L8").Value = "MEDELLÍN"
$ws.Range("M8").Value = "053324501001900119"
$ws.Range("N8").Value = "Medellín"
$ws.Range("O8").Value = "AEJDO"
$ws.Range("P8").Value = "Habilitación Viviendas Metrosur"
$ws.Range("Q8").Value = "Urbano"
$ws.Range("R8").Value = 5
$ws.Range("S8").Value = 45946.6625
$ws.Range("T8").Value = "16 días 15:54"
$ws.Range("U8").Value = "VENCIDO"
$ws.Range("V8").Value = "VENCIDO"
$ws.Range("W8").Value = "CERRADO"
$ws.Range("X8").Value = "Ejecutado en Campo"

# Row 9
$ws.Range("A9").Value = "23332144"
$ws.Range("B9").Value = "ENENOR"
$ws.Range("C9").Value = "NUEVO"
$ws.Range("D9").Value = "ENENOR"
$ws.Range("E9").Value = "20/01/2025 09:02"
$ws.Range("F9").Value = 45930.31180555555
$ws.Range("G9").Value = 8355854
$ws.Range("H9").Value = "JUAN CARLOS LOPEZ MOLINA"
$ws.Range("I9").Value = 9876543
$ws.Range("J9").Value = "SIN DATOS"
$ws.Range("K9").Value = "RURAL_140004950000000049_PROV.PARCELACION VOLTA HO"
$ws.Range("L9").Value = "ENVIGADO"
$ws.Range("M9").Value = "140004950000000049"
$ws.Range("N9").Value = "Medellín"
$ws.Range("O9").Value = "ALEGN"
$ws.Range("P9").Value = "Revisor Instalac. Regadas Oriente"
$ws.Range("Q9").Value = "Rural"
$ws.Range("R9").Value = 10
$ws.Range("S9").Value = 45945.31180555555
$ws.Range("T9").Value = "22 días 07:29"
$ws.Range("U9").Value = "VENCIDO"
$ws.Range("V9").Value = "VENCIDO"
$ws.Range("W9").Value = "CERRADO"
$ws.Range("X9").Value = "Ejecutado en Campo"

# Row 10
$ws.Range("A10").Value = "23499958"
$ws.Range("B10").Value = "ENEMRT"
$ws.Range("C10").Value = "NUEVO"
$ws.Range("D10").Value = "ENEMVI"
$ws.Range("E10").Value = "29/07/2025 14:48"
$ws.Range("F10").Value = 45894.37777777778
$ws.Range("G10").Value = 71590457
$ws.Range("H10").Value = "CARLOS ALBERTO JARAMILLO MESA"
$ws.Range("I10").Value = "SIN DATOS"
$ws.Range("J10").Value = 3113651861
$ws.Range("K10").Value = "RURAL_147014002000000000_147014002000000000"
$ws.Range("L10").Value = "MEDELLÍN"
$ws.Range("M10").Value = "147014002000000000"
$ws.Range("N10").Value = "Medellín"
$ws.Range("O10").Value = "AMRTR"
$ws.Range("P10").Value = "MET-RMRT-Francisco J Dominguez"
$ws.Range("Q10").Value = "Rural"
$ws.Range("R10").Value = 10
$ws.Range("S10").Value = 45908.37777777778
$ws.Range("T10").Value = "48 días 09:04"
$ws.Range("U10").Value = "VENCIDO"
$ws.Range("V10").Value = "VENCIDO"
$ws.Range("W10").Value = "CERRADO"
$ws.Range("X10").Value = "Ejecutado en Campo"

# Restore the Normal style on A5:A10 / M5:M10 so only the number format used
# for text entry is dropped (no stray "@" text style is left on the cell) --
# the stored cell content remains text either way.
$ws.Range("A5:A10").Style = "Normal"
$ws.Range("M5:M10").Style = "Normal"
